# Programming Support Agreed-Upon Standards - apply commit edits
# 1) Remove the stray _GoBack bookmark sitting in the empty Heading1
#    paragraph near the top of the document.
# 2) Append a "New objects" sub-section at the end of the document
#    (a bullet, a sub-bullet whose run ends with a new _GoBack bookmark,
#    and three blank/indented paragraphs).

$d = $word.ActiveDocument

# --- Step 1: drop the old _GoBack bookmark -------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# --- Step 2: insert the new paragraphs before the trailing empty ---------
# paragraph that sits just before the sectPr.
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $pLast.Range
$r.Collapse(1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr></w:pPr><w:r><w:t>New objects</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="11"/></w:numPr></w:pPr><w:r><w:t>All objects must have a completed and approved Design Specification form</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + `
  '<w:p><w:pPr><w:ind w:left="720"/></w:pPr></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="2160"/></w:pPr></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="2160"/></w:pPr></w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
